$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook appears to track a rolling window of quarterly error-metric rows
# (Q1..Q9) against ifoCAST. This edit adds the next quarter's full-series
# evaluation: every existing row's data shifts up by one (row N takes what was
# row N+1), and a brand new final row (row 11) is populated with the newest
# quarter's error stats, while the N counters in column G each decrease by 1.

$data = @{
    2  = @{ B = -0.1143386297721907;  C = 0.6889840097079422;  D = 0.9692665404913205;  E = 0.9845133521142924;  F = 1.006200714633619;  G = 18 }
    3  = @{ B = -0.03231471279853988; C = 0.7695157101427289;  D = 0.9504010021523444;  E = 0.9748851225412891;  F = 1.004336377456358;   G = 17 }
    4  = @{ B = 0.05671622629984657;  C = 0.6321910596607474;  D = 0.8518300650585396;  E = 0.9229464042177854;  F = 0.951413452211165;   G = 16 }
    5  = @{ B = 0.085928964333323;    C = 0.5612671955956002;  D = 0.5873845025919633;  E = 0.7664101399328973;  F = 0.7883078955770884;  G = 15 }
    6  = @{ B = 0.1382138394267656;   C = 0.7476276879241753;  D = 1.097975515794379;   E = 1.047843268716452;   F = 1.077897322974905;   G = 14 }
    7  = @{ B = 0.1215206328360093;   C = 0.7819526689445239;  D = 1.205421270947194;   E = 1.097916786895616;   F = 1.135726716235405;   G = 13 }
    8  = @{ B = 0.2028461735751207;   C = 0.8472735835465873;  D = 1.158512388347609;   E = 1.076342133500128;   F = 1.10405820904526;    G = 12 }
    9  = @{ B = 0.239907996146195;    C = 0.6651045147971374;  D = 0.5928596743010739;  E = 0.769973814035954;   F = 0.7673553351966808;  G = 11 }
    10 = @{ B = 0.2225839694609408;   C = 0.8306889768957216;  D = 1.245987321698922;   E = 1.116238021973325;   F = 1.152988242321246;   G = 10 }
    11 = @{ B = 0.2059614158583246;   C = 1.002733503272531;   D = 1.92968265145822;    E = 1.389130178010045;   F = 1.457110278930643;   G = 9 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
}
